# Auto-generated script to apply Twintania_Profits.xlsx market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 129.66667
$ws.Range("I2").Value = 145.6
$ws.Range("K2").Value = 145.6
$ws.Range("M2").Value = -32.59999999999999
$ws.Range("H6").Value = 3081.0715
$ws.Range("I6").Value = 363.57144
$ws.Range("J6").Value = 5798.5713
$ws.Range("K6").Value = 1090.71432
$ws.Range("L6").Value = 17395.7139
$ws.Range("M6").Value = -978.71432
$ws.Range("N6").Value = -17619.7139
$ws.Range("H33").Value = 6761757
$ws.Range("I33").Value = 27027028
$ws.Range("J33").Value = 6666.3335
$ws.Range("K33").Value = 27027028
$ws.Range("L33").Value = 6666.3335
$ws.Range("M33").Value = -27026799
$ws.Range("N33").Value = -7124.3335
$ws.Range("H40").Value = 2139.8
$ws.Range("I40").Value = 2103.0303
$ws.Range("K40").Value = 2103.0303
$ws.Range("M40").Value = -1928.0303
$ws.Range("H64").Value = 6398.2856
$ws.Range("I64").Value = 3600
$ws.Range("K64").Value = 3600
$ws.Range("M64").Value = -3352
$ws.Range("H67").Value = 6398.2856
$ws.Range("I67").Value = 3600
$ws.Range("K67").Value = 3600
$ws.Range("M67").Value = -2742
$ws.Range("H70").Value = 40650.375
$ws.Range("I70").Value = 1733.3334
$ws.Range("J70").Value = 64000.6
$ws.Range("K70").Value = 5200.0002
$ws.Range("L70").Value = 192001.8
$ws.Range("M70").Value = -4930.0002
$ws.Range("N70").Value = -192541.8
$ws.Range("H73").Value = 40650.375
$ws.Range("I73").Value = 1733.3334
$ws.Range("J73").Value = 64000.6
$ws.Range("K73").Value = 5200.0002
$ws.Range("L73").Value = 192001.8
$ws.Range("M73").Value = -4264.0002
$ws.Range("N73").Value = -193873.8
$ws.Range("H80").Value = 257169.77
$ws.Range("I80").Value = 553
$ws.Range("J80").Value = 556556
$ws.Range("K80").Value = 1659
$ws.Range("L80").Value = 1669668
$ws.Range("M80").Value = -661
$ws.Range("N80").Value = -1671664
$ws.Range("H83").Value = 257169.77
$ws.Range("I83").Value = 553
$ws.Range("J83").Value = 556556
$ws.Range("K83").Value = 4977
$ws.Range("L83").Value = 5009004
$ws.Range("M83").Value = 15
$ws.Range("N83").Value = -5018988
$ws.Range("H98").Value = 2679.88
$ws.Range("I98").Value = 2743.087
$ws.Range("J98").Value = 1953
$ws.Range("K98").Value = 2743.087
$ws.Range("L98").Value = 1953
$ws.Range("M98").Value = -1245.087
$ws.Range("N98").Value = -4949
$ws.Range("H100").Value = 45109.61
$ws.Range("I100").Value = 46705.5
$ws.Range("K100").Value = 46705.5
$ws.Range("M100").Value = -46164.5
$ws.Range("H122").Value = 2679.88
$ws.Range("I122").Value = 2743.087
$ws.Range("J122").Value = 1953
$ws.Range("K122").Value = 8229.261
$ws.Range("L122").Value = 5859
$ws.Range("M122").Value = -5779.261
$ws.Range("N122").Value = -10759
$ws.Range("H132").Value = 2385.5833
$ws.Range("I132").Value = 2334.6897
$ws.Range("K132").Value = 7004.0691
$ws.Range("M132").Value = -4474.0691
$ws.Range("H137").Value = 8147.609
$ws.Range("I137").Value = 3583.6843
$ws.Range("J137").Value = 13742.097
$ws.Range("K137").Value = 10751.0529
$ws.Range("L137").Value = 41226.291
$ws.Range("M137").Value = -8201.052899999999
$ws.Range("N137").Value = -46326.291
$ws.Range("H138").Value = 3330.95
$ws.Range("I138").Value = 1495.4445
$ws.Range("J138").Value = 3512.4834
$ws.Range("K138").Value = 4486.333500000001
$ws.Range("L138").Value = 10537.4502
$ws.Range("M138").Value = 653.6664999999994
$ws.Range("N138").Value = -20817.4502
$ws.Range("H141").Value = 2852.0715
$ws.Range("I141").Value = 2852.0715
$ws.Range("K141").Value = 8556.2145
$ws.Range("M141").Value = -3376.2145
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 30830.666
$ws.Range("I28").Value = 30830.666
$ws.Range("K28").Value = 30830.666
$ws.Range("M28").Value = -30638.666
$ws.Range("H32").Value = 6965.933
$ws.Range("I32").Value = 6729.5957
$ws.Range("K32").Value = 6729.5957
$ws.Range("M32").Value = -6442.5957
$ws.Range("H34").Value = 28012.5
$ws.Range("I34").Value = 28012.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 28012.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -27741.5
$ws.Range("N34").ClearContents()
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 20000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -20626
$ws.Range("H61").Value = 9049.538
$ws.Range("I61").Value = 6077.1665
$ws.Range("K61").Value = 6077.1665
$ws.Range("M61").Value = -5865.1665
$ws.Range("H74").Value = 9969.957
$ws.Range("I74").Value = 8709.755999999999
$ws.Range("K74").Value = 8709.755999999999
$ws.Range("M74").Value = -7835.755999999999
$ws.Range("H77").Value = 9969.957
$ws.Range("I77").Value = 8709.755999999999
$ws.Range("K77").Value = 43548.78
$ws.Range("M77").Value = -39180.78
$ws.Range("H99").Value = 30830.666
$ws.Range("I99").Value = 30830.666
$ws.Range("K99").Value = 30830.666
$ws.Range("M99").Value = -27835.666
$ws.Range("H102").Value = 3372.95
$ws.Range("I102").Value = 1113.0714
$ws.Range("K102").Value = 1113.0714
$ws.Range("M102").Value = 508.9286
$ws.Range("H110").Value = 3585.2666
$ws.Range("I110").Value = 3179.3
$ws.Range("J110").Value = 4397.2
$ws.Range("K110").Value = 3179.3
$ws.Range("L110").Value = 4397.2
$ws.Range("M110").Value = -1134.3
$ws.Range("N110").Value = -8487.200000000001
$ws.Range("H122").Value = 3814.1365
$ws.Range("I122").Value = 3244.85
$ws.Range("K122").Value = 9734.549999999999
$ws.Range("M122").Value = -7284.549999999999
$ws.Range("H132").Value = 2593.3147
$ws.Range("I132").Value = 2705.7954
$ws.Range("J132").Value = 2098.4
$ws.Range("K132").Value = 8117.3862
$ws.Range("L132").Value = 6295.200000000001
$ws.Range("M132").Value = -5587.3862
$ws.Range("N132").Value = -11355.2
$ws.Range("H136").Value = 9049.538
$ws.Range("I136").Value = 6077.1665
$ws.Range("K136").Value = 18231.4995
$ws.Range("M136").Value = -15681.4995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 877.1177
$ws.Range("I80").Value = 561.1429000000001
$ws.Range("J80").Value = 1098.3
$ws.Range("K80").Value = 561.1429000000001
$ws.Range("L80").Value = 1098.3
$ws.Range("M80").Value = 436.8570999999999
$ws.Range("N80").Value = -3094.3
$ws.Range("H83").Value = 877.1177
$ws.Range("I83").Value = 561.1429000000001
$ws.Range("J83").Value = 1098.3
$ws.Range("K83").Value = 2805.7145
$ws.Range("L83").Value = 5491.5
$ws.Range("M83").Value = 2186.2855
$ws.Range("N83").Value = -15475.5
$ws.Range("H86").Value = 336751.6
$ws.Range("I86").Value = 834753.25
$ws.Range("J86").Value = 4750.5
$ws.Range("K86").Value = 834753.25
$ws.Range("L86").Value = 4750.5
$ws.Range("M86").Value = -833630.25
$ws.Range("N86").Value = -6996.5
$ws.Range("H89").Value = 336751.6
$ws.Range("I89").Value = 834753.25
$ws.Range("J89").Value = 4750.5
$ws.Range("K89").Value = 4173766.25
$ws.Range("L89").Value = 23752.5
$ws.Range("M89").Value = -4168150.25
$ws.Range("N89").Value = -34984.5
$ws.Range("H105").Value = 3809.1614
$ws.Range("I105").Value = 3474.625
$ws.Range("K105").Value = 3474.625
$ws.Range("M105").Value = -1727.625
$ws.Range("H107").Value = 4336
$ws.Range("I107").Value = 4308
$ws.Range("K107").Value = 4308
$ws.Range("M107").Value = -2388
$ws.Range("H132").Value = 84963.336
$ws.Range("J132").Value = 84963.336
$ws.Range("L132").Value = 84963.336
$ws.Range("N132").Value = -95083.336
$ws.Range("H134").Value = 13130.615
$ws.Range("I134").Value = 7588.8125
$ws.Range("K134").Value = 22766.4375
$ws.Range("M134").Value = -20231.4375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 3000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = -2623
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2789
$ws.Range("H62").Value = 256003.25
$ws.Range("I62").Value = 338002.34
$ws.Range("J62").Value = 10006
$ws.Range("K62").Value = 338002.34
$ws.Range("L62").Value = 10006
$ws.Range("M62").Value = -337378.34
$ws.Range("N62").Value = -11254
$ws.Range("H65").Value = 256003.25
$ws.Range("I65").Value = 338002.34
$ws.Range("J65").Value = 10006
$ws.Range("K65").Value = 1690011.7
$ws.Range("L65").Value = 50030
$ws.Range("M65").Value = -1686891.7
$ws.Range("N65").Value = -56270
$ws.Range("H86").Value = 3286.0967
$ws.Range("I86").Value = 2747.158
$ws.Range("J86").Value = 4139.4165
$ws.Range("K86").Value = 2747.158
$ws.Range("L86").Value = 4139.4165
$ws.Range("M86").Value = -1624.158
$ws.Range("N86").Value = -6385.4165
$ws.Range("H89").Value = 3286.0967
$ws.Range("I89").Value = 2747.158
$ws.Range("J89").Value = 4139.4165
$ws.Range("K89").Value = 13735.79
$ws.Range("L89").Value = 20697.0825
$ws.Range("M89").Value = -8119.789999999999
$ws.Range("N89").Value = -31929.0825
$ws.Range("H107").Value = 1339.7333
$ws.Range("I107").Value = 1267.3846
$ws.Range("J107").Value = 1810
$ws.Range("K107").Value = 1267.3846
$ws.Range("L107").Value = 1810
$ws.Range("M107").Value = 652.6153999999999
$ws.Range("N107").Value = -5650
$ws.Range("H132").Value = 34400.848
$ws.Range("I132").Value = 28207.387
$ws.Range("J132").Value = 43543.57
$ws.Range("K132").Value = 84622.16099999999
$ws.Range("L132").Value = 130630.71
$ws.Range("M132").Value = -82092.16099999999
$ws.Range("N132").Value = -135690.71
$ws.Range("H134").Value = 2739.8809
$ws.Range("I134").Value = 1306.1786
$ws.Range("J134").Value = 5607.2856
$ws.Range("K134").Value = 3918.5358
$ws.Range("L134").Value = 16821.8568
$ws.Range("M134").Value = -1383.5358
$ws.Range("N134").Value = -21891.8568
$ws.Range("H141").Value = 460198.78
$ws.Range("J141").Value = 651643.1
$ws.Range("L141").Value = 651643.1
$ws.Range("N141").Value = -662003.1
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 137.33333
$ws.Range("J23").Value = 134.8
$ws.Range("L23").Value = 404.4
$ws.Range("N23").Value = -874.4000000000001
$ws.Range("H26").Value = 56.666668
$ws.Range("I26").Value = 90.888885
$ws.Range("J26").Value = 36.133335
$ws.Range("K26").Value = 272.666655
$ws.Range("L26").Value = 108.400005
$ws.Range("M26").Value = 15.33334500000001
$ws.Range("N26").Value = -684.400005
$ws.Range("H38").Value = 3212.4707
$ws.Range("I38").Value = 488.33334
$ws.Range("J38").Value = 6277.125
$ws.Range("K38").Value = 1465.00002
$ws.Range("L38").Value = 18831.375
$ws.Range("M38").Value = -1118.00002
$ws.Range("N38").Value = -19525.375
$ws.Range("H55").Value = 5036
$ws.Range("J55").Value = 9922.5
$ws.Range("L55").Value = 29767.5
$ws.Range("N55").Value = -30121.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H107").Value = 1527.25
$ws.Range("J107").Value = 1087.8
$ws.Range("L107").Value = 3263.4
$ws.Range("N107").Value = -7103.4
$ws.Range("H123").Value = 2828
$ws.Range("I123").Value = 1893.6
$ws.Range("J123").Value = 7500
$ws.Range("K123").Value = 5680.799999999999
$ws.Range("L123").Value = 22500
$ws.Range("M123").Value = -3230.799999999999
$ws.Range("N123").Value = -27400
$ws.Range("H126").Value = 14900
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 14900
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 44700
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -54580
$ws.Range("H133").Value = 4399.6
$ws.Range("I133").Value = 4399.6
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 13198.8
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -8138.800000000001
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 463
$ws.Range("I134").Value = 463
$ws.Range("K134").Value = 1389
$ws.Range("M134").Value = 3681
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 1650
$ws.Range("J25").Value = 1650
$ws.Range("L25").Value = 1650
$ws.Range("N25").Value = -2708
$ws.Range("H33").Value = 23237
$ws.Range("I33").Value = 8986
$ws.Range("J33").Value = 26799.75
$ws.Range("K33").Value = 8986
$ws.Range("L33").Value = 26799.75
$ws.Range("M33").Value = -8734
$ws.Range("N33").Value = -27303.75
$ws.Range("H70").Value = 10213
$ws.Range("I70").Value = 8346.333000000001
$ws.Range("J70").Value = 12313
$ws.Range("K70").Value = 8346.333000000001
$ws.Range("L70").Value = 12313
$ws.Range("M70").Value = -8076.333000000001
$ws.Range("N70").Value = -12853
$ws.Range("H73").Value = 10213
$ws.Range("I73").Value = 8346.333000000001
$ws.Range("J73").Value = 12313
$ws.Range("K73").Value = 8346.333000000001
$ws.Range("L73").Value = 12313
$ws.Range("M73").Value = -7410.333000000001
$ws.Range("N73").Value = -14185
$ws.Range("H113").Value = 186107.9
$ws.Range("I113").Value = 227031
$ws.Range("J113").Value = 1954
$ws.Range("K113").Value = 227031
$ws.Range("L113").Value = 1954
$ws.Range("M113").Value = -224861
$ws.Range("N113").Value = -6294
$ws.Range("H122").Value = 3715.0688
$ws.Range("I122").Value = 4253.4707
$ws.Range("K122").Value = 12760.4121
$ws.Range("M122").Value = -10310.4121
$ws.Range("H126").Value = 6993.8
$ws.Range("I126").Value = 4989.6665
$ws.Range("K126").Value = 14968.9995
$ws.Range("M126").Value = -12498.9995
$ws.Range("H132").Value = 8513.958000000001
$ws.Range("I132").Value = 7516.75
$ws.Range("K132").Value = 22550.25
$ws.Range("M132").Value = -20020.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6035.2383
$ws.Range("I7").Value = 4476.8
$ws.Range("K7").Value = 4476.8
$ws.Range("M7").Value = -4364.8
$ws.Range("H40").Value = 8294.053
$ws.Range("I40").Value = 8652.235000000001
$ws.Range("K40").Value = 8652.235000000001
$ws.Range("M40").Value = -8516.235000000001
$ws.Range("H93").Value = 6560.4688
$ws.Range("I93").Value = 5969.8076
$ws.Range("K93").Value = 5969.8076
$ws.Range("M93").Value = -4721.8076
$ws.Range("H122").Value = 6041.375
$ws.Range("I122").Value = 4870.533
$ws.Range("J122").Value = 7992.778
$ws.Range("K122").Value = 14611.599
$ws.Range("L122").Value = 23978.334
$ws.Range("M122").Value = -12161.599
$ws.Range("N122").Value = -28878.334
$ws.Range("H126").Value = 6035.2383
$ws.Range("I126").Value = 4476.8
$ws.Range("K126").Value = 13430.4
$ws.Range("M126").Value = -10960.4
$ws.Range("H132").Value = 5695.8667
$ws.Range("I132").Value = 5337.6313
$ws.Range("K132").Value = 16012.8939
$ws.Range("M132").Value = -13482.8939
$ws.Range("H136").Value = 6110.81
$ws.Range("I136").Value = 5672.969
$ws.Range("K136").Value = 17018.907
$ws.Range("M136").Value = -14468.907
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 84145.664
$ws.Range("J34").Value = 84737.5
$ws.Range("L34").Value = 84737.5
$ws.Range("N34").Value = -85143.5
$ws.Range("H42").Value = 67047.336
$ws.Range("I42").Value = 41144
$ws.Range("J42").Value = 79999
$ws.Range("K42").Value = 41144
$ws.Range("L42").Value = 79999
$ws.Range("M42").Value = -40766
$ws.Range("N42").Value = -80755
$ws.Range("H43").Value = 66996.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 66996.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 66996.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -67294.5
$ws.Range("H62").Value = 9677.357
$ws.Range("J62").Value = 8708.272000000001
$ws.Range("L62").Value = 8708.272000000001
$ws.Range("N62").Value = -9956.272000000001
$ws.Range("H65").Value = 9677.357
$ws.Range("J65").Value = 8708.272000000001
$ws.Range("L65").Value = 43541.36
$ws.Range("N65").Value = -49781.36
$ws.Range("H75").Value = 49990.5
$ws.Range("I75").Value = 49990.5
$ws.Range("K75").Value = 49990.5
$ws.Range("M75").Value = -49054.5
$ws.Range("H78").Value = 49990.5
$ws.Range("I78").Value = 49990.5
$ws.Range("K78").Value = 149971.5
$ws.Range("M78").Value = -145291.5
$ws.Range("H107").Value = 2965.111
$ws.Range("I107").Value = 2177.2
$ws.Range("K107").Value = 6531.599999999999
$ws.Range("M107").Value = -4611.599999999999
$ws.Range("H113").Value = 725.5294
$ws.Range("I113").Value = 737.5
$ws.Range("K113").Value = 2212.5
$ws.Range("M113").Value = -42.5
$ws.Range("H126").Value = 9368.895
$ws.Range("I126").Value = 6235.0938
$ws.Range("K126").Value = 18705.2814
$ws.Range("M126").Value = -16235.2814
$ws.Range("H132").Value = 117250.06
$ws.Range("I132").Value = 163514.45
$ws.Range("K132").Value = 490543.35
$ws.Range("M132").Value = -488013.35
$ws.Range("H136").Value = 10348244
$ws.Range("I136").Value = 11541411
$ws.Range("K136").Value = 34624233
$ws.Range("M136").Value = -34621683
